$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 165, shifting existing rows 165-218 down to 166-219
$ws.Rows("165:165").Insert()

# Fill the new row 165 with its data
$ws.Range("A165").Value = 4
$ws.Range("B165").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C165").Value = "Los Lagos"
$ws.Range("D165").Value = 44627
$ws.Range("E165").Value = 10
$ws.Range("F165").Value = 100112044
$ws.Range("G165").Value = "Perejil"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 80
$ws.Range("K165").Value = 5000
$ws.Range("L165").Value = 5000
$ws.Range("M165").Value = 5000
$ws.Range("N165").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O165").Value = "Región Metropolitana"
$ws.Range("P165").Value = 1667
$ws.Range("Q165").Value = 3
$ws.Range("R165").Value = "Hortaliza"
